$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"25.93000000000061"
$ws.Range("H2").Value = [double]"1.110223024625157e-16"
$ws.Range("I2").Value = [double]"1.110223024625157e-16"
$ws.Range("L2").Value = [double]"49.41438068660503"
$ws.Range("M2").Value = "[41.80582809805051, 57.022933275159545]"
$ws.Range("N2").Value = [double]"0"
$ws.Range("O2").Value = [double]"0"
$ws.Range("P2").Value = [double]"1.591237119836271"
$ws.Range("Q2").Value = "[1.415131825941347, 1.7673424137311953]"
$ws.Range("T2").Value = [double]"56.78536377946965"
$ws.Range("U2").Value = "[51.491603877724145, 62.07912368121515]"
$ws.Range("X2").Value = [double]"19.3631431431436"
$ws.Range("Y2").Value = [double]"18.63637637637682"
$ws.Range("Z2").Value = [double]"20.08990990991039"
$ws.Range("F3").Value = [double]"25.93000000000061"
$ws.Range("L3").Value = [double]"49.45537084780003"
$ws.Range("M3").Value = "[40.882252557439216, 58.028489138160836]"
$ws.Range("N3").Value = [double]"3.996802888650564e-15"
$ws.Range("O3").Value = [double]"3.996802888650564e-15"
$ws.Range("P3").Value = [double]"1.113237036407194"
$ws.Range("Q3").Value = "[0.9245527929483472, 1.3019212798660416]"
$ws.Range("R3").Value = [double]"1.77635683940025e-15"
$ws.Range("S3").Value = [double]"1.77635683940025e-15"
$ws.Range("T3").Value = [double]"54.20307528575066"
$ws.Range("U3").Value = "[49.50150810563903, 58.904642465862295]"
$ws.Range("X3").Value = [double]"21.3357957957963"
$ws.Range("Y3").Value = [double]"20.5571171171176"
$ws.Range("Z3").Value = [double]"22.114474474475"
$ws.Range("F4").Value = [double]"25.93000000000061"
$ws.Range("L4").Value = [double]"47.76388271258746"
$ws.Range("M4").Value = "[38.42024822896587, 57.10751719620905]"
$ws.Range("N4").Value = [double]"2.069455717901292e-13"
$ws.Range("O4").Value = [double]"2.069455717901292e-13"
$ws.Range("P4").Value = [double]"0.6478159025420398"
$ws.Range("Q4").Value = "[0.4465527095192714, 0.8490790955648082]"
$ws.Range("R4").Value = [double]"5.977741701634898e-08"
$ws.Range("S4").Value = [double]"5.977741701634898e-08"
$ws.Range("T4").Value = [double]"52.55116765716807"
$ws.Range("U4").Value = "[47.65792306599356, 57.44441224834257]"
$ws.Range("X4").Value = [double]"23.25653653653708"
$ws.Range("Y4").Value = [double]"22.42594594594648"
$ws.Range("Z4").Value = [double]"24.08712712712769"
$ws.Range("F5").Value = [double]"25.93000000000061"
$ws.Range("L5").Value = [double]"50.9850779587748"
$ws.Range("M5").Value = "[43.802517785162756, 58.167638132386834]"
$ws.Range("N5").Value = [double]"0"
$ws.Range("O5").Value = [double]"0"
$ws.Range("P5").Value = [double]"0.2327105669325773"
$ws.Range("Q5").Value = "[0.09434212172942313, 0.37107901213573147]"
$ws.Range("R5").Value = [double]"0.001474933046643034"
$ws.Range("S5").Value = [double]"0.001474933046643034"
$ws.Range("T5").Value = [double]"51.34935474202927"
$ws.Range("U5").Value = "[47.247824195584585, 55.450885288473955]"
$ws.Range("X5").Value = [double]"24.96962962963022"
$ws.Range("Y5").Value = [double]"24.39859859859918"
$ws.Range("Z5").Value = [double]"25.54066066066126"
$ws.Range("F6").Value = [double]"25.93000000000061"
$ws.Range("H6").Value = [double]"7.771561172376096e-16"
$ws.Range("I6").Value = [double]"7.771561172376096e-16"
$ws.Range("L6").Value = [double]"48.80274068658083"
$ws.Range("M6").Value = "[39.539013660920574, 58.06646771224108]"
$ws.Range("N6").Value = [double]"7.860379014346108e-14"
$ws.Range("O6").Value = [double]"7.860379014346108e-14"
$ws.Range("P6").Value = [double]"-0.2641579408423853"
$ws.Range("Q6").Value = "[-0.4654211338651546, -0.06289474781961601]"
$ws.Range("R6").Value = [double]"0.01125252462954762"
$ws.Range("S6").Value = [double]"0.01125252462954762"
$ws.Range("T6").Value = [double]"51.82571507844427"
$ws.Range("U6").Value = "[46.518043003548705, 57.13338715333983]"
$ws.Range("X6").Value = [double]"1.090150150150176"
$ws.Range("Y6").Value = [double]"0.259559559559567"
$ws.Range("Z6").Value = [double]"1.920740740740785"
$ws.Range("F7").Value = [double]"25.93000000000061"
$ws.Range("H7").Value = [double]"1.110223024625157e-16"
$ws.Range("I7").Value = [double]"1.110223024625157e-16"
$ws.Range("L7").Value = [double]"49.03181919604843"
$ws.Range("M7").Value = "[40.2675970475373, 57.79604134455956]"
$ws.Range("N7").Value = [double]"1.088018564132653e-14"
$ws.Range("O7").Value = [double]"1.088018564132653e-14"
$ws.Range("P7").Value = [double]"-0.754736973835386"
$ws.Range("Q7").Value = "[-0.9434212172942322, -0.5660527303765397]"
$ws.Range("R7").Value = [double]"2.853726144280699e-10"
$ws.Range("S7").Value = [double]"2.853726144280699e-10"
$ws.Range("T7").Value = [double]"49.00111459843973"
$ws.Range("U7").Value = "[44.3543231588791, 53.64790603800037]"
$ws.Range("X7").Value = [double]"3.11471471471479"
$ws.Range("Y7").Value = [double]"2.336036036036093"
$ws.Range("Z7").Value = [double]"3.893393393393486"
$ws.Range("F8").Value = [double]"25.93000000000061"
$ws.Range("H8").Value = [double]"4.440892098500626e-16"
$ws.Range("I8").Value = [double]"4.440892098500626e-16"
$ws.Range("L8").Value = [double]"42.43925765896378"
$ws.Range("M8").Value = "[34.25760751485052, 50.62090780307703]"
$ws.Range("N8").Value = [double]"1.296740492762183e-13"
$ws.Range("O8").Value = [double]"1.296740492762183e-13"
$ws.Range("P8").Value = [double]"-1.182421259008771"
$ws.Range("Q8").Value = "[-1.3962634015954647, -0.9685791164220783]"
$ws.Range("R8").Value = [double]"1.598721155460225e-14"
$ws.Range("S8").Value = [double]"1.598721155460225e-14"
$ws.Range("T8").Value = [double]"52.34712383438866"
$ws.Range("U8").Value = "[47.54267895106629, 57.15156871771103]"
$ws.Range("X8").Value = [double]"4.879719719719837"
$ws.Range("Y8").Value = [double]"3.997217217217311"
$ws.Range("Z8").Value = [double]"5.762222222222363"
$ws.Range("F9").Value = [double]"22.79000000000012"
$ws.Range("H9").Value = [double]"1.110223024625157e-16"
$ws.Range("I9").Value = [double]"1.110223024625157e-16"
$ws.Range("L9").Value = [double]"50.65397432344014"
$ws.Range("M9").Value = "[43.21495074702016, 58.092997899860116]"
$ws.Range("N9").Value = [double]"0"
$ws.Range("O9").Value = [double]"0"
$ws.Range("P9").Value = [double]"-1.735895039821387"
$ws.Range("Q9").Value = "[-1.8868424345884645, -1.58494764505431]"
$ws.Range("T9").Value = [double]"52.10518716326305"
$ws.Range("U9").Value = "[47.55864012576009, 56.651734200766]"
$ws.Range("X9").Value = [double]"6.29633633633637"
$ws.Range("Y9").Value = [double]"5.748828828828859"
$ws.Range("Z9").Value = [double]"6.843843843843882"
$ws.Range("F10").Value = [double]"22.79000000000012"
$ws.Range("H10").Value = [double]"6.561418075534675e-14"
$ws.Range("I10").Value = [double]"6.561418075534675e-14"
$ws.Range("L10").Value = [double]"46.32895757491627"
$ws.Range("M10").Value = "[36.64447368261111, 56.01344146722144]"
$ws.Range("N10").Value = [double]"1.644462344074782e-12"
$ws.Range("O10").Value = [double]"1.644462344074782e-12"
$ws.Range("Q10").Value = "[-1.62268449374608, -1.1698423094448476]"
$ws.Range("R10").Value = [double]"4.440892098500626e-16"
$ws.Range("S10").Value = [double]"4.440892098500626e-16"
$ws.Range("T10").Value = [double]"51.6130271735116"
$ws.Range("U10").Value = "[45.90657533495265, 57.31947901207056]"
$ws.Range("X10").Value = [double]"5.064444444444472"
$ws.Range("Y10").Value = [double]"4.243183183183205"
$ws.Range("Z10").Value = [double]"5.88570570570574"
$ws.Range("F11").Value = [double]"22.79000000000012"
$ws.Range("L11").Value = [double]"48.98748807105589"
$ws.Range("M11").Value = "[39.58728198389846, 58.38769415821333]"
$ws.Range("N11").Value = [double]"1.116884362772907e-13"
$ws.Range("O11").Value = [double]"1.116884362772907e-13"
$ws.Range("P11").Value = [double]"-1.081789662497386"
$ws.Range("Q11").Value = "[-1.2830528555201557, -0.8805264694746171]"
$ws.Range("R11").Value = [double]"4.085620730620576e-14"
$ws.Range("S11").Value = [double]"4.085620730620576e-14"
$ws.Range("T11").Value = [double]"52.79726904329942"
$ws.Range("U11").Value = "[47.781433031605296, 57.81310505499355]"
$ws.Range("X11").Value = [double]"3.923803803803825"
$ws.Range("Y11").Value = [double]"3.193793793793811"
$ws.Range("Z11").Value = [double]"4.653813813813838"
